$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row after the existing data (row 30 is the last used row -> new rows 31, 32)
$lastRow = $ws.Cells.Item(1, 1).End(4).Row
$row1 = $lastRow + 1
$row2 = $lastRow + 2

# --- Row 1 (Paris Lowres): columns A-F ---
$ws.Cells.Item($row1, 1).Value = "Paris Lowres"
$ws.Cells.Item($row1, 2).Value = "SSW Office Assistant"
$ws.Cells.Item($row1, 3).Value = "Recruitment"
$ws.Cells.Item($row1, 4).Value = "Scrum"

# --- Row 2 (Christian Morford-Waite): columns A-F ---
$ws.Cells.Item($row2, 1).Value = "Christian Morford-Waite"
$ws.Cells.Item($row2, 2).Value = "SSW Software Architect"
$ws.Cells.Item($row2, 3).Value = ".NET Core"
$ws.Cells.Item($row2, 4).Value = "Azure DevOps"
$ws.Cells.Item($row2, 5).Value = "Power BI"

# --- Profile column (G) for both rows ---
$ws.Cells.Item($row1, 7).Value = "Paris is a copy editor at university and working for SSW while finishing her degree. Paris makes the Melbourne office run smoothly!"
$ws.Cells.Item($row2, 7).Value = "​​​​​​​​​​Christian is a Software Developer with over 5 years of experience working with a wide range of technologies and Azure resources.`nSpecialising in backend .NET development and API integration. He enjoys working on process improvement and automation through PowerShell scripting and Azure Build Pipelines.`n"

# --- Remaining columns (H-M) for both rows ---
$ws.Cells.Item($row1, 8).Value = 205
$ws.Cells.Item($row1, 9).Value = "ssw_tv"
$ws.Cells.Item($row1, 10).Value = $true
$ws.Cells.Item($row1, 11).Value = "N"
$ws.Cells.Item($row1, 12).Value = "N"
$ws.Cells.Item($row1, 13).Value = "N"

$ws.Cells.Item($row2, 8).Value = 205
$ws.Cells.Item($row2, 9).Value = "ssw_tv"
$ws.Cells.Item($row2, 10).Value = $true
$ws.Cells.Item($row2, 11).Value = "N"
$ws.Cells.Item($row2, 12).Value = "N"
$ws.Cells.Item($row2, 13).Value = "N"

# Match row height/formatting of the previous data row for the two new rows
$ws.Rows.Item($row1).RowHeight = $ws.Rows.Item($lastRow).RowHeight
$ws.Rows.Item($row2).RowHeight = $ws.Rows.Item($lastRow).RowHeight

$srcRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 13))
$dstRange1 = $ws.Range($ws.Cells.Item($row1, 1), $ws.Cells.Item($row1, 13))
$dstRange2 = $ws.Range($ws.Cells.Item($row2, 1), $ws.Cells.Item($row2, 13))
$srcRange.Copy()
$dstRange1.PasteSpecial(-4122)
$dstRange2.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the final selection left by the editing session
[void]$ws.Range("G30").Select()
